$wb = $excel.ActiveWorkbook

$wsInstructions = $wb.Worksheets.Item("Instructions")
$wsDataset      = $wb.Worksheets.Item("Dataset")
$wsTerminology  = $wb.Worksheets.Item("Terminology")

# ---------------------------------------------------------------------------
# 1) Instructions sheet: reword the three "TODO" bullet lines that described
#    the Ka/Kd/KD columns and now describe "Standard deviation" columns.
#    The sheet is protected, so unprotect -> edit -> re-protect.
# ---------------------------------------------------------------------------
$wsInstructions.Unprotect()
$wsInstructions.Range("A10").Value = "- Standard deviation in M^-1s^-1: TODO"
$wsInstructions.Range("A12").Value = "- Standard deviation in 1/s: TODO"
$wsInstructions.Range("A14").Value = "- Standard deviation in nM: TODO"
$wsInstructions.Protect()

# ---------------------------------------------------------------------------
# 2) Dataset sheet: rename the duplicated header cells (E1/G1/I1) to
#    "Standard deviation ..." and widen those columns; add a list data
#    validation on J2:J100 backed by the Terminology sheet.
# ---------------------------------------------------------------------------
$wsDataset.Range("E1").Value = "Standard deviation in M^-1s^-1"
$wsDataset.Range("G1").Value = "Standard deviation in 1/s"
$wsDataset.Range("I1").Value = "Standard deviation in nM"

# Column widths (OOXML stored width = ColumnWidth + 5/6)
$wsDataset.Columns.Item(5).ColumnWidth = 30 - 5/6   # E -> 30
$wsDataset.Columns.Item(7).ColumnWidth = 25 - 5/6   # G -> 25
$wsDataset.Columns.Item(9).ColumnWidth = 24 - 5/6   # I -> 24

$wsDataset.Range("J2:J100").Validation.Add(3, 1, 1, "=Terminology!`$A`$2:`$A`$4")

# ---------------------------------------------------------------------------
# 3) Terminology sheet: populate the lookup list used by the new validation
#    (header "Qualitiative measure" in bold + positive/negative/unknown),
#    freeze the header row and set the column width. Sheet is protected too.
# ---------------------------------------------------------------------------
$wsTerminology.Unprotect()

$wsInstructions.Range("A1").Copy()
$wsTerminology.Range("A1").PasteSpecial(-4122)

$wsTerminology.Range("A1").Value = "Qualitiative measure"
$wsTerminology.Range("A2").Value = "positive"
$wsTerminology.Range("A3").Value = "negative"
$wsTerminology.Range("A4").Value = "unknown"

$wsTerminology.Columns.Item(1).ColumnWidth = 20 - 5/6   # A -> 20

$wsTerminology.Activate()
$wsTerminology.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$wsTerminology.Range("A1").Select()

$wsTerminology.Protect()

# Restore the sheet that was active before these edits (Dataset), so the
# workbook-level "active tab" stays as close as possible to the original.
$wsDataset.Activate()

